$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Resistor" to "Resistors" in Part Number (A2) and Library Ref (B2)
$ws.Range("A2").Value = "Resistors"
$ws.Range("B2").Value = "Resistors"

# Update the selected cell to C7
$ws.Range("C7").Select()
